# Fruta / hortaliza, semanal
# Re-order the weekly price rows (2-8) for Granada at Lo Valledor so the
# dates/data line up with the latest week's observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current (pre-edit) values for columns D, K, L, M, N, O, P, Q, R, S, T
# across rows 2-8, keyed by row number.
$rows = 2..8
$data = @{}
foreach ($r in $rows) {
    $data[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        R = $ws.Cells.Item($r, 18).Value2
        S = $ws.Cells.Item($r, 19).Value2
        T = $ws.Cells.Item($r, 20).Value2
    }
}

# Mapping of new-row -> old-row (where the data now sitting in new-row came from).
$mapping = @{
    2 = 5
    3 = 6
    4 = 7
    5 = 8
    6 = 2
    7 = 4
    8 = 3
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $src = $data[$oldRow]

    $ws.Cells.Item($newRow, 4).Value2 = $src.D
    $ws.Cells.Item($newRow, 11).Value2 = $src.K
    $ws.Cells.Item($newRow, 12).Value2 = $src.L
    $ws.Cells.Item($newRow, 13).Value2 = $src.M
    $ws.Cells.Item($newRow, 14).Value2 = $src.N
    $ws.Cells.Item($newRow, 15).Value2 = $src.O
    $ws.Cells.Item($newRow, 16).Value2 = $src.P
    $ws.Cells.Item($newRow, 17).Value2 = $src.Q
    $ws.Cells.Item($newRow, 18).Value2 = $src.R
    $ws.Cells.Item($newRow, 19).Value2 = $src.S
    $ws.Cells.Item($newRow, 20).Value2 = $src.T
}
